# Apply updated odds/stat values to the "Jogos da Semana" FlashScore sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("H3").Value  = 5.75
$ws.Range("AH3").Value = 26
$ws.Range("AM3").Value = 101
$ws.Range("AP3").Value = 23
$ws.Range("AY3").Value = 67

# Row 4 updates
$ws.Range("J4").Value  = 2.38
$ws.Range("L4").Value  = 5
$ws.Range("M4").Value  = 1.06
$ws.Range("N4").Value  = 10
$ws.Range("O4").Value  = 1.3
$ws.Range("P4").Value  = 3.4
$ws.Range("Q4").Value  = 1.98
$ws.Range("R4").Value  = 1.83
$ws.Range("Z4").Value  = 13
$ws.Range("AB4").Value = 29
$ws.Range("AI4").Value = 23
$ws.Range("AL4").Value = 41
$ws.Range("AW4").Value = 6.5
$ws.Range("AX4").Value = 26
